$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2278
    $ws.Range("F3").Value = 1729
    $ws.Range("F4").Value = 337
    $ws.Range("F5").Value = 1093
    $ws.Range("F6").Value = 856
    $ws.Range("F8").Value = 5847
}

$wb.Save()
